$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.835940330277484, 3.230376277753256, 0.2926989868425439),
    @(8.475344398840846, 6.984832544430162, 0.1046356402809086),
    @(8.631195580249331, 3.866418119113546, 0.2429153396576904),
    @(8.997468905687182, 11.3864114120676, 0.3149368528947745),
    @(11.21206103842722, 4.735509233805131, 0.3987609710512867),
    @(14.57935990438859, 4.665168682740525, 0.2555202946935587),
    @(14.87497347225885, 4.915375584805513, 0.3419357241511858),
    @(15.50837519354749, 5.416161294987182, 0.08521242169328522),
    @(16.77482751291398, 3.667717151358299, 0.3215056188415228),
    @(20.63139664295238, 3.740597772363054, 0.122830270515992),
    @(21.15192989415073, 10.21363055831688, 0.3534461488825369),
    @(21.33388482737919, 4.369581885073277, 0.1619881376898507),
    @(24.80017405620643, 11.7129360566776, 0.5424233266480257),
    @(25.56522779962446, 6.326122436069289, 0.455298509422811),
    @(26.43341987300064, 2.623697903572965, 0.277239177483033),
    @(31.09041086878132, 5.363855985624959, 0.2054052315177273),
    @(31.59327420912789, 11.02752812786597, 0.2734367037740995),
    @(32.68921439478268, 6.733329698549494, 0.3226102352803845),
    @(33.62909348313258, 5.174453786977046, 0.217977132522423),
    @(35.17494858068797, 2.379255306931858, 0.2518266547266559),
    @(38.72201377512928, 10.38845423934469, 0.4406309967226275),
    @(44.47617890710871, 5.543518137976344, 0.4276398685039972),
    @(45.59664606451351, 1.525279926936533, 0.1938432043782048),
    @(59.55584277247033, 8.4652273662902, 0.3073265388170409),
    @(66.83395720587539, 3.726437243746259, 0.1746197229544027),
    @(67.39615534098206, 5.481992472932705, 0.2064033648773851),
    @(71.00993007968667, 8.807975441925445, 0.2941528578173399),
    @(87.9829781893532, 4.584239571246671, 0.3987171273803903),
    @(90.3029266862508, 5.491126954588989, 0.2703690812815168),
    @(90.49540817524218, 2.797966978456842, 0.1782442861792367),
    @(93.61222139446551, 8.022678061788417, 0.4295377321902156),
    @(94.09459006280413, 4.07555108224582, 0.2847670942756964),
    @(98.53939003495401, 1.449928269533457, 0.2154738764362104),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
}

$ws.Range("A35:D48").EntireRow.Delete()
